$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "최종점수" (Final Score, column K) for rows 2-6
$ws.Range("K2").Value = 58.5
$ws.Range("K3").Value = 57.3
$ws.Range("K4").Value = 50.5
$ws.Range("K5").Value = 48.3
$ws.Range("K6").Value = 45.3

# Update "MACRO_SCORE" (column N) for rows 2-6
$ws.Range("N2").Value = 51.15965480231979
$ws.Range("N3").Value = 51.15965480231979
$ws.Range("N4").Value = 51.15965480231979
$ws.Range("N5").Value = 51.15965480231979
$ws.Range("N6").Value = 51.15965480231979
